$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Line 1 Hourly Production (column B) - fill in the previously-blank 4pm/4:30pm/5pm hours
$ws.Range("B9").Value = 114
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100

# Line 2 Hourly Production (column G) - fill in the previously-blank 4pm/4:30pm/5pm hours
$ws.Range("G9").Value = 202
$ws.Range("G10").Value = 203
$ws.Range("G11").Value = 203

# Line 3 Hourly Production (column L) - fill in all hourly production figures for the day
$ws.Range("L2").Value = 100
$ws.Range("L3").Value = 203
$ws.Range("L4").Value = 4
$ws.Range("L5").Value = 390
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 29
$ws.Range("L8").Value = 289
$ws.Range("L9").Value = 298
$ws.Range("L10").Value = 38
$ws.Range("L11").Value = 239

# Move the active selection to B11 to match the saved view state
$ws.Range("B11").Select()
